$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Altınyıldız Classics"
$ws.Range("B1").Value = "249,99 TL"
$ws.Range("D1").Value = "https://cdn.dsmcdn.com/mnresize/420/620/ty80/product/media/images/20210308/18/69852728/137294414/1/1_org_zoom.jpg"
$ws.Range("E1").Value = "['https://cdn.dsmcdn.com/mnresize/420/620/ty81/product/media/images/20210308/18/69852728/137294414/2/2_org_zoom.jpg']"

$ws.Range("A2").Value = "Happiness İstanbul"
$ws.Range("B2").Value = "239,48 TL"
$ws.Range("D2").Value = "https://cdn.dsmcdn.com/mnresize/420/620/ty496/product/media/images/20220801/13/153049070/533113392/2/2_org_zoom.jpg"
$ws.Range("E2").Value = "['https://cdn.dsmcdn.com/mnresize/420/620/ty497/product/media/images/20220801/13/153049070/533113392/3/3_org_zoom.jpg', 'https://cdn.dsmcdn.com/mnresize/420/620/ty497/product/media/images/20220801/13/153049070/533113392/4/4_org_zoom.jpg', 'https://cdn.dsmcdn.com/mnresize/420/620/ty497/product/media/images/20220801/13/153049070/533113392/5/5_org_zoom.jpg', 'https://cdn.dsmcdn.com/mnresize/420/620/ty495/product/media/images/20220801/13/153049070/533113392/6/6_org_zoom.jpg']"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""

$ws.Range("A3").Value = "Default Title"
$ws.Range("B3").Value = "163,29 TL"
$ws.Range("C3").Value = "165,64 TL"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "['https://cdn.dsmcdn.com/mnresize/1200/1800/ty1162/product/media/images/prod/SPM/PIM/20240206/15/bf8ceb40-539b-3e2d-9aa0-555191680f22/1_org_zoom.jpg', 'https://cdn.dsmcdn.com/ty1162/product/media/images/prod/SPM/PIM/20240206/15/bf8ceb40-539b-3e2d-9aa0-555191680f22/1_org_zoom.jpg']"
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = ""

$ws.Range("A4").Value = "Default Title"
$ws.Range("B4").Value = "247,87 TL"
$ws.Range("E4").Value = "['https://cdn.dsmcdn.com/mnresize/1200/1800/ty999/product/media/images/prod/PIM/20230918/11/e24750c5-08e8-4457-9e97-07d1bf7a4181/1_org_zoom.jpg', 'https://cdn.dsmcdn.com/ty999/product/media/images/prod/PIM/20230918/11/e24750c5-08e8-4457-9e97-07d1bf7a4181/1_org_zoom.jpg']"

$ws.Range("A5").Value = "Default Title"
$ws.Range("B5").Value = "349 TL"
$ws.Range("D5").Value = "https://cdn.dsmcdn.com/mnresize/420/620/ty1143/product/media/images/prod/SPM/PIM/20240122/19/d6b9e8ea-eadb-3352-b11c-b793f27e5dca/1_org_zoom.jpg"
$ws.Range("E5").Value = "['https://cdn.dsmcdn.com/mnresize/420/620/ty1144/product/media/images/prod/SPM/PIM/20240122/19/5ecad0a9-639b-377d-ab88-5e043c1da4bf/1_org_zoom.jpg', 'https://cdn.dsmcdn.com/mnresize/420/620/ty1144/product/media/images/prod/SPM/PIM/20240122/19/182ea6e2-b62a-36e0-8c14-3336328a1bfe/1_org_zoom.jpg', 'https://cdn.dsmcdn.com/mnresize/420/620/ty1144/product/media/images/prod/SPM/PIM/20240122/19/0175491e-2c51-3c1a-a074-7427bc490c03/1_org_zoom.jpg', 'https://cdn.dsmcdn.com/mnresize/420/620/ty1145/product/media/images/prod/SPM/PIM/20240122/19/1f66bb55-7ad0-309e-ba61-24f499f4e7c8/1_org_zoom.jpg']"

$ws.Range("A6").Value = "Twisted Minds"
$ws.Range("B6").Value = "1.205 TL"
$ws.Range("D6").Value = "https://cdn.dsmcdn.com/mnresize/420/620/ty1006/product/media/images/prod/SPM/PIM/20230929/09/d5f4cec3-e06c-3ae7-b3d5-d61fbc79db36/1_org_zoom.jpg"
$ws.Range("E6").Value = "['https://cdn.dsmcdn.com/mnresize/420/620/ty1005/product/media/images/prod/SPM/PIM/20230929/09/26f15a9e-199c-3f2a-8cce-d071aed38376/1_org_zoom.jpg', 'https://cdn.dsmcdn.com/mnresize/420/620/ty1006/product/media/images/prod/SPM/PIM/20230929/09/82502fea-32ab-35fa-8784-68c79fb5552f/1_org_zoom.jpg', 'https://cdn.dsmcdn.com/mnresize/420/620/ty1006/product/media/images/prod/SPM/PIM/20230929/09/4810b651-27e1-3663-92b2-f966cbe62f48/1_org_zoom.jpg', 'https://cdn.dsmcdn.com/mnresize/420/620/ty1005/product/media/images/prod/SPM/PIM/20230929/09/4312aaf8-ba56-3379-b367-ba3a0950e8bb/1_org_zoom.jpg']"

$ws.Range("A7").Value = "Default Title"
$ws.Range("B7").Value = "2.519 TL"
$ws.Range("C7").Value = "2.919 TL"
$ws.Range("E7").Value = "['https://cdn.dsmcdn.com/mnresize/1200/1800/ty1026/product/media/images/prod/SPM/PIM/20231030/13/2e233745-406f-348a-a33f-6086a5d6a8ad/1_org_zoom.jpg', 'https://cdn.dsmcdn.com/ty1026/product/media/images/prod/SPM/PIM/20231030/13/2e233745-406f-348a-a33f-6086a5d6a8ad/1_org_zoom.jpg']"
